$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows to append (date serial, B, C, D) starting at row 344,
# continuing the series through 23 Aug 2021 (commit: "aggiornamento a l 23 agosto 2021")
$data = @(
    @(44418, 6, 27, 157.5722206011089),
    @(44419, 0, 26, 151.7362124306974),
    @(44420, 3, 27, 157.5722206011089),
    @(44421, 5, 28, 163.4082287715203),
    @(44422, 5, 29, 169.2442369419317),
    @(44423, 2, 26, 151.7362124306974),
    @(44424, 1, 22, 128.3921797490516),
    @(44425, 2, 18, 105.0481470674059),
    @(44426, 0, 18, 105.0481470674059),
    @(44427, 5, 20, 116.7201634082288),
    @(44428, 2, 17, 99.21213889699445),
    @(44429, 1, 13, 75.8681062153487),
    @(44430, 1, 12, 70.03209804493727),
    @(44431, 3, 14, 81.70411438576014)
)

$startRow = 344
$templateRow = $startRow - 1

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $startRow + $i
    $row = $data[$i]

    # Copy cell A formatting (date style) from the row above, same as the
    # existing data rows, then overwrite with the new value.
    $ws.Cells.Item($templateRow, 1).Copy($ws.Cells.Item($r, 1))
    $ws.Cells.Item($r, 1).Value = $row[0]

    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}
